$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row3 = @(26011018018, 13, 2, 0, 0, 1, 2, 1, 0, 0, 0, 0, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 1, 1, 3, 0.75, 4, 0)
for ($c = 1; $c -le $row3.Length; $c++) {
  $ws.Cells.Item(3, $c).Value = $row3[$c - 1]
}

$row4 = @(26021035008, 14, 2, 0, 0, 1, 3, 0, 1, 0, 0, 1, 0, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 0, 1, 2, 0.5, 4, 1)
for ($c = 1; $c -le $row4.Length; $c++) {
  $ws.Cells.Item(4, $c).Value = $row4[$c - 1]
}

$row5 = @(26021035009, 13, 1, 0, 0, 1, 3, 0, 0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 1, 2, 2, 0.7, 2, 0)
for ($c = 1; $c -le $row5.Length; $c++) {
  $ws.Cells.Item(5, $c).Value = $row5[$c - 1]
}

$row6 = @(26021035013, 14, 2, 0, 0, 1, 3, 0, 0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 0, 1, 3, 0.77, 2, 0)
for ($c = 1; $c -le $row6.Length; $c++) {
  $ws.Cells.Item(6, $c).Value = $row6[$c - 1]
}

$row7 = @(26021035024, 15, 2, 0, 0, 1, 4, 0, 0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 0, 1, 0, 0, 1, 1, 0, 1, 3, 0.69, 3, 0)
for ($c = 1; $c -le $row7.Length; $c++) {
  $ws.Cells.Item(7, $c).Value = $row7[$c - 1]
}

$row8 = @(26031048033, 12, 1, 0, 1, 1, 1, 0, 1, 0, 0, 0, 1, 0, 0, 0, 1, 1, 1, 1, 0, 0, 1, 1, 0, 1, 2, 0.59, 4, 0)
for ($c = 1; $c -le $row8.Length; $c++) {
  $ws.Cells.Item(8, $c).Value = $row8[$c - 1]
}

$row9 = @(26022036034, 15, 1, 0, 0, 4, 3, 0, 0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 0, 0, 3, 0.63, 3, 1)
for ($c = 1; $c -le $row9.Length; $c++) {
  $ws.Cells.Item(9, $c).Value = $row9[$c - 1]
}

$row10 = @(26022036035, 14, 1, 0, 0, 4, 4, 0, 0, 1, 0, 1, 0, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 0, 2, 2, 0.64, 4, 0)
for ($c = 1; $c -le $row10.Length; $c++) {
  $ws.Cells.Item(10, $c).Value = $row10[$c - 1]
}

$row11 = @(26022036036, 14, 1, 0, 0, 1, 3, 0, 0, 1, 0, 0, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 0, 0, 3, 0.06, 3, 0)
for ($c = 1; $c -le $row11.Length; $c++) {
  $ws.Cells.Item(11, $c).Value = $row11[$c - 1]
}

$row12 = @(26022036037, 14, 2, 0, 0, 1, 3, 0, 0, 1, 0, 1, 0, 0, 0, 0, 1, 1, 0, 0, 0, 0, 1, 1, 0, 2, 3, 0.19, 2, 0)
for ($c = 1; $c -le $row12.Length; $c++) {
  $ws.Cells.Item(12, $c).Value = $row12[$c - 1]
}

$ws.AutoFilterMode = $false
$ws.Range("F9").Select()
